$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add footnote markers: "Random Forest" -> "Random Forest*" and
# "RF Lower Cutoff" -> "*RF Lower Cutoff (decreasing the probability from 50% to 30%)"
$ws.Range("A4").Value = "Random Forest*"
$ws.Range("A5").Value = "*RF Lower Cutoff (decreasing the probability from 50% to 30%)"

# Widen column A a bit and grow row 5's height to fit the longer note text
# (42.15 -> stored width snaps to a clean 43 after Excel's pixel rounding)
$ws.Columns("A").ColumnWidth = 42.15
$ws.Rows("5").RowHeight = 63

# Move the active selection
$ws.Range("B10").Select()
